$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measurement values in row 2
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = 101
$ws.Range("C2").Value = 3

# Move the active selection to C2 (previously B3)
$ws.Range("C2").Select()
